$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127; existing rows 127-214 shift down to 128-215.
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row 127 with the new record.
$ws.Cells.Item(127, 1).Value = 4
$ws.Cells.Item(127, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(127, 3).Value = "Los Lagos"
$ws.Cells.Item(127, 4).Value = 44897
$ws.Cells.Item(127, 5).Value = 10
$ws.Cells.Item(127, 6).Value = "Fruta"
$ws.Cells.Item(127, 7).Value = 100103
$ws.Cells.Item(127, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(127, 9).Value = 100103004
$ws.Cells.Item(127, 10).Value = "Durazno"
$ws.Cells.Item(127, 11).Value = "Florida King"
$ws.Cells.Item(127, 12).Value = "Primera"
$ws.Cells.Item(127, 13).Value = 600
$ws.Cells.Item(127, 14).Value = 23000
$ws.Cells.Item(127, 15).Value = 24000
$ws.Cells.Item(127, 16).Value = 23500
$ws.Cells.Item(127, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(127, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(127, 19).Value = 1679
$ws.Cells.Item(127, 20).Value = 14
